# Workbook produced two small edits on "sample1.xlsx":
#  - Sheet1 gained a 4th data row (F4:H4), re-using existing shared strings
#    str7/str8/str9 (the same values already used on Sheet2 row 2).
#  - The active selection moved: Sheet1's cursor ended on H14 (no longer the
#    selected tab) and Sheet2 became the active tab with C1 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Make sure Sheet1 is active while we edit/select on it.
$ws1.Select()

$ws1.Range("F4").Value = "str7"
$ws1.Range("G4").Value = "str8"
$ws1.Range("H4").Value = "str9"

# Leave Sheet1's selection parked on H14 (below/right of the used range).
$ws1.Range("H14").Select()

# Switch to Sheet2 and select C1 there - this is the sheet/cell that ends
# up active when the workbook is saved.
$ws2.Select()
$ws2.Range("C1").Select()
